$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update metric values
$ws.Range("B3").Value = 9068982222090.75
$ws.Range("C3").Value = 13751630899183.12
$ws.Range("D3").Value = 11238905383591.78

# Row 4: rename model and update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 4326911561266.629
$ws.Range("C4").Value = 4326911561266.627
$ws.Range("D4").Value = 4326911561266.628

# Row 5: rename model and update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 443601960344741.6
$ws.Range("C5").Value = 405173521700405.2
$ws.Range("D5").Value = 476228807997274.3
